$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "-0.25**"
$ws.Range("C3").Value = "-0.353***"
$ws.Range("C4").Value = "4.659***"
$ws.Range("C5").Value = 0.52
